# Insert a new weekly price record as row 102, pushing the existing rows
# (old rows 102-186) down by one (new rows 103-187).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(102).Insert()

$ws.Range("A102").Value = 3
$ws.Range("B102").Value = "Femacal de La Calera"
$ws.Range("C102").Value = "Coquimbo"
$ws.Range("D102").Value = 44447
$ws.Range("E102").Value = 5
$ws.Range("F102").Value = 100114013
$ws.Range("G102").Value = "Zanahoria"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 410
$ws.Range("K102").Value = 5000
$ws.Range("L102").Value = 5500
$ws.Range("M102").Value = 5305
$ws.Range("N102").Value = "$/saco 20 kilos"
$ws.Range("O102").Value = "Provincia de Quillota"
$ws.Range("P102").Value = 265
$ws.Range("Q102").Value = 20
$ws.Range("R102").Value = "Hortaliza"
